$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financial Forecast")
$ws.Range("A26:AM26").Insert(-4121)
Write-Host "done"
